$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "774×9=" "113×8="
Replace-Text "518×3=" "429×9="
Replace-Text "688×3=" "290×8="
Replace-Text "588×4=" "263×4="
Replace-Text "810×5=" "111×6="
Replace-Text "932×3=" "603×7="
Replace-Text "334×4=" "707×3="
Replace-Text "662×3=" "107×2="
Replace-Text "542×3=" "786×4="
Replace-Text "150×2=" "996×2="
Replace-Text "437×6=" "218×5="
Replace-Text "799×4=" "212×9="
Replace-Text "854×4=" "323×4="
Replace-Text "819×9=" "830×5="
Replace-Text "317×5=" "920×3="
Replace-Text "917×2=" "831×5="
Replace-Text "521×7=" "848×9="
Replace-Text "825×5=" "467×2="
Replace-Text "792×9=" "892×4="
Replace-Text "745×8=" "946×9="
Replace-Text "644×5=" "358×6="
Replace-Text "659×7=" "644×4="
Replace-Text "964×2=" "133×2="
Replace-Text "437×7=" "584×7="
Replace-Text "984×8=" "741×2="
